$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 157
$data = @(
    @(155, 43938, "LAS TUNAS", 0, 0, 0, 0),
    @(156, 43938, "HOLGUÍN", 15, 0, 1, 0),
    @(157, 43938, "GRANMA", 0, 0, 0, 0),
    @(158, 43938, "SANTIAGO", 73, 0, 2, 1),
    @(159, 43938, "GUANTÁNAMO", 0, 0, 0, 0)
)

$ws.Range("A156:G156").Copy() | Out-Null
$ws.Range("A157:G161").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$ws.Range("F161").Select() | Out-Null
